$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "id" column (A) was missing "pub-18" (it jumped straight from pub-17 to
# pub-19), while a stray duplicate "pub-73" sat unused at the bottom of the
# shared-string table. Shifting every id in rows 19-73 down by one fills the
# pub-18 gap and removes the need for the stray pub-73 entry.
for ($row = 19; $row -le 73; $row++) {
    $ws.Cells.Item($row, 1).Value = "pub-" + ($row - 1)
}

# Restore the cursor/selection state left behind by the editor.
$ws.Range("B76").Select() | Out-Null
